$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns stay text, matching the source data format
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.679.64"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "1.804.62"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "313.81"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.5407"
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D8").Value = "0.3789"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D9").Value = "0.07537"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "42.65"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").Value = "1.120"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "21.01"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "6.185"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "7.397"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("D16").Value = "1.795.48"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "90.71"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "0.00001066"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "0.06446"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "17.27"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "5.930"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "28.671.54"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("D24").Value = "11.20"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "2.112"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").Value = "161.15"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("D27").Value = "20.53"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "2.383"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "2.007.85"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "123.60"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "1.112"
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("D32").Value = "0.1042"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "5.687"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("D34").Value = "3.697"
$ws.Range("D35").Value = "0.2264"
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("D36").Value = "0.06496"
$ws.Range("E36").Value = "  +7.36%  "
$ws.Range("D37").Value = "8.958"
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").Value = "5.058"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "11.34"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("D41").Value = "0.6265"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").Value = "1.203"
$ws.Range("E42").Value = "  +4.50%  "
$ws.Range("D43").Value = "0.9985"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "1.395"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "13.38"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "0.5892"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "3.671"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "126.27"
$ws.Range("E48").Value = "  +3.38%  "
$ws.Range("D49").Value = "1.963"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").Value = "1.156"
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("D51").Value = "0.06894"
$ws.Range("E51").Value = "  +1.50%  "
